# MirrorMe/Arguments.xlsx — add "Templates" and "Variables" sheets
# (accompanies issue #266)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Tweak the F15 phrasing text on "Arguments" (drop the trailing space)
#    and move the selection the author ended up leaving on that sheet.
# ---------------------------------------------------------------------
$wsArguments = $wb.Worksheets.Item("Arguments")
$wsArguments.Range("F15").Value = "De werknemer heeft aanspraak op een verhoging."

# ---------------------------------------------------------------------
# 2. "Statements" selection moved to A3:A10
# ---------------------------------------------------------------------
$wsStatements = $wb.Worksheets.Item("Statements")
$wsStatements.Activate()
$wsStatements.Range("A3:A10").Select()

# ---------------------------------------------------------------------
# 3. Add the "Templates" sheet (after the last existing sheet)
# ---------------------------------------------------------------------
$wsTemplates = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTemplates.Name = "Templates"

$wsTemplates.Range("A1").Value = "[Template]"
$wsTemplates.Range("B1").Value = "phrasing"
$wsTemplates.Range("A2").Value = "Template"
$wsTemplates.Range("B2").Value = "Text"

$wsTemplates.Range("A3").Value = "t1"
$wsTemplates.Range("B3").Value = "prefix 1"
$wsTemplates.Range("A4").Value = "t2"
$wsTemplates.Range("B4").Value = "prefix 2"

$wsTemplates.Range("A17").Value = "[Sequence]"
$wsTemplates.Range("B17").Value = "head"
$wsTemplates.Range("A18").Value = "Sequence"
$wsTemplates.Range("B18").Value = "Item"

$wsTemplates.Range("A19").Value = "s1"
$wsTemplates.Range("A20").Value = "s2"
$wsTemplates.Range("A21").Value = "s3"
$wsTemplates.Range("A22").Value = "s4"
$wsTemplates.Range("A23").Value = "s5"
$wsTemplates.Range("A24").Value = "s6"
$wsTemplates.Range("A25").Value = "s7"
$wsTemplates.Range("A26").Value = "s8"

$wsTemplates.Range("C3").Select()

# ---------------------------------------------------------------------
# 4. Add the "Variables" sheet (after "Templates")
# ---------------------------------------------------------------------
$wsVariables = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsVariables.Name = "Variables"

$wsVariables.Range("A1").Value = "[Variable]"
$wsVariables.Range("B1").Value = "type"
$wsVariables.Range("A2").Value = "Variable"
$wsVariables.Range("B2").Value = "Concept"

$wsVariables.Range("A3").Value = "v1"
$wsVariables.Range("B3").Value = "Persoon"
$wsVariables.Range("A4").Value = "v2"
$wsVariables.Range("B4").Value = "Persoon"
$wsVariables.Range("A5").Value = "v3"
$wsVariables.Range("B5").Value = "Persoon"

$wsVariables.Range("A14").Value = "[Item]"
$wsVariables.Range("B14").Value = "succ"
$wsVariables.Range("C14").Value = "seq"
$wsVariables.Range("A15").Value = "Item"
$wsVariables.Range("B15").Value = "Item"
$wsVariables.Range("C15").Value = "Sequence"

$wsVariables.Range("A16").Value = "v1"
$wsVariables.Range("B16").Value = "v2"
$wsVariables.Range("C16").Value = "t1"
$wsVariables.Range("A17").Value = "v2"
$wsVariables.Range("B17").Value = "v3"
$wsVariables.Range("C17").Value = "t1"
$wsVariables.Range("A18").Value = "v3"
$wsVariables.Range("C18").Value = "t1"

$wsVariables.Range("B14").Select()

# ---------------------------------------------------------------------
# 5. Leave "Arguments" as the active sheet/tab, as in the original file.
# ---------------------------------------------------------------------
$wsArguments.Activate()
$wsArguments.Range("E22").Select()
